$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.190.37'
$ws.Range("E2").Value = '  -0.86%  '

$ws.Range("D3").Value = '2.243.49'
$ws.Range("E3").Value = '  -1.74%  '

$ws.Range("E4").Value = '  +0.04%  '

$c = $ws.Range("D5")
$c.Value = "'246.93"
$c.Style = "Normal"
$ws.Range("E5").Value = '  -1.38%  '

$c = $ws.Range("D6")
$c.Value = "'0.630"
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.40%  '

$c = $ws.Range("D7")
$c.Value = "'76.14"
$c.Style = "Normal"
$ws.Range("E7").Value = '  +5.83%  '

$ws.Range("E8").Value = '  +0.02%  '

$c = $ws.Range("D9")
$c.Value = "'0.632"
$c.Style = "Normal"
$ws.Range("E9").Value = '  -0.04%  '

$c = $ws.Range("D10")
$c.Value = "'40.44"
$c.Style = "Normal"
$ws.Range("E10").Value = '  +4.72%  '

$ws.Range("E11").Value = '  -1.74%  '

$c = $ws.Range("D12")
$c.Value = "'7.23"
$c.Style = "Normal"
$ws.Range("E12").Value = '  -1.14%  '

$ws.Range("E13").Value = '  -2.00%  '

$ws.Range("D14").Value = '2.578.57'
$ws.Range("E14").Value = '  -1.87%  '

$c = $ws.Range("D15")
$c.Value = "'14.89"
$c.Style = "Normal"
$ws.Range("E15").Value = '  -1.12%  '

$c = $ws.Range("D16")
$c.Value = "'0.861"
$c.Style = "Normal"
$ws.Range("E16").Value = '  -2.14%  '

$ws.Range("D17").Value = '2.239.15'
$ws.Range("E17").Value = '  -1.73%  '

$ws.Range("D18").Value = '42.156.19'
$ws.Range("E18").Value = '  -0.93%  '

$ws.Range("D19").Value = '0.0₃0978'
$ws.Range("E19").Value = '  -1.71%  '

$ws.Range("E20").Value = '  -1.88%  '

$c = $ws.Range("D21")
$c.Value = "'71.51"
$c.Style = "Normal"
$ws.Range("E21").Value = '  -1.34%  '

$ws.Range("B22").Value = 'BitcoinCash'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Range("D22")
$c.Value = "'231.75"
$c.Style = "Normal"
$ws.Range("E22").Value = '  -1.38%  '

$ws.Range("B23").Value = 'ImmutableX'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range("D23")
$c.Value = "'2.21"
$c.Style = "Normal"
$ws.Range("E23").Value = '  -4.08%  '

$ws.Range("E24").Value = '  -0.04%  '

$ws.Range("E25").Value = '  -5.15%  '

$c = $ws.Range("D26")
$c.Value = "'11.16"
$c.Style = "Normal"
$ws.Range("E26").Value = '  -3.79%  '

$ws.Range("E27").Value = '  -4.81%  '

$c = $ws.Range("D28")
$c.Value = "'7.07"
$c.Style = "Normal"
$ws.Range("E28").Value = '  +10.74%  '

$ws.Range("E29").Value = '  -1.40%  '

$c = $ws.Range("D30")
$c.Value = "'168.60"
$c.Style = "Normal"
$ws.Range("E30").Value = '  +0.71%  '

$c = $ws.Range("D31")
$c.Value = "'20.53"
$c.Style = "Normal"
$ws.Range("E31").Value = '  -2.47%  '

$c = $ws.Range("D32")
$c.Value = "'0.0855"
$c.Style = "Normal"
$ws.Range("E32").Value = '  +6.72%  '

$c = $ws.Range("D33")
$c.Value = "'32.70"
$c.Style = "Normal"
$ws.Range("E33").Value = '  +2.96%  '

$ws.Range("E34").Value = '  -6.77%  '

$ws.Range("E35").Value = '  -0.09%  '

$ws.Range("E36").Value = '  -4.64%  '

$c = $ws.Range("D37")
$c.Value = "'4.80"
$c.Style = "Normal"
$ws.Range("E37").Value = '  +1.61%  '

$c = $ws.Range("D38")
$c.Value = "'0.0297"
$c.Style = "Normal"
$ws.Range("E38").Value = '  -3.00%  '

$c = $ws.Range("D39")
$c.Value = "'13.38"
$c.Style = "Normal"
$ws.Range("E39").Value = '  -3.74%  '

$c = $ws.Range("D40")
$c.Value = "'5.94"
$c.Style = "Normal"
$ws.Range("E40").Value = '  -0.59%  '

$ws.Range("E41").Value = '  -5.97%  '

$c = $ws.Range("D42")
$c.Value = "'117.23"
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.Value = "'0.203"
$c.Style = "Normal"
$ws.Range("E43").Value = '  -4.53%  '

$c = $ws.Range("D44")
$c.Value = "'60.10"
$c.Style = "Normal"
$ws.Range("E44").Value = '  -2.82%  '

$c = $ws.Range("D45")
$c.Value = "'8.72"
$c.Style = "Normal"
$ws.Range("E45").Value = '  -5.78%  '

$ws.Range("E46").Value = '  -2.60%  '

$c = $ws.Range("D47")
$c.Value = "'0.998"
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.25%  '

$ws.Range("E48").Value = '  -4.00%  '

$ws.Range("E49").Value = '  -1.04%  '

$c = $ws.Range("D50")
$c.Value = "'4.24"
$c.Style = "Normal"
$ws.Range("E50").Value = '  -13.11%  '

$c = $ws.Range("D51")
$c.Value = "'2.26"
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.71%  '
